$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2158.875
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2158.875
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6476.625
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6812.625

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H64").Value = 4999.5
$ws.Range("I64").Value = 4999.5
$ws.Range("K64").Value = 4999.5
$ws.Range("M64").Value = -4751.5

$ws.Range("H67").Value = 4999.5
$ws.Range("I67").Value = 4999.5
$ws.Range("K67").Value = 4999.5
$ws.Range("M67").Value = -4141.5

$ws.Range("H70").Value = 1648.75
$ws.Range("I70").Value = 1631.6666
$ws.Range("J70").Value = 1700
$ws.Range("K70").Value = 4894.9998
$ws.Range("L70").Value = 5100
$ws.Range("M70").Value = -4624.9998
$ws.Range("N70").Value = -5640

$ws.Range("H73").Value = 1648.75
$ws.Range("I73").Value = 1631.6666
$ws.Range("J73").Value = 1700
$ws.Range("K73").Value = 4894.9998
$ws.Range("L73").Value = 5100
$ws.Range("M73").Value = -3958.9998
$ws.Range("N73").Value = -6972

$ws.Range("H100").Value = 6463.5713
$ws.Range("I100").Value = 2649
$ws.Range("J100").Value = 16000
$ws.Range("K100").Value = 2649
$ws.Range("L100").Value = 16000
$ws.Range("M100").Value = -2108
$ws.Range("N100").Value = -17082

$ws.Range("H132").Value = 1308.2609
$ws.Range("I132").Value = 1205
$ws.Range("K132").Value = 3615
$ws.Range("M132").Value = -1085

$ws.Range("H135").Value = 1018.8947
$ws.Range("I135").Value = 964.3889
$ws.Range("K135").Value = 8679.500100000001
$ws.Range("M135").Value = -6144.500100000001

$ws.Range("H137").Value = 4032.3333
$ws.Range("I137").Value = 4032.3333
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 12096.9999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -9546.999899999999
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 4441.7856
$ws.Range("I138").Value = 2865.625
$ws.Range("J138").Value = 4812.647
$ws.Range("K138").Value = 8596.875
$ws.Range("L138").Value = 14437.941
$ws.Range("M138").Value = -3456.875
$ws.Range("N138").Value = -24717.941

$ws.Range("H141").Value = 3555.2222
$ws.Range("I141").Value = 3119.9333
$ws.Range("J141").Value = 5731.6665
$ws.Range("K141").Value = 9359.7999
$ws.Range("L141").Value = 17194.9995
$ws.Range("M141").Value = -4179.7999
$ws.Range("N141").Value = -27554.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3773.2454
$ws.Range("I32").Value = 3595.64
$ws.Range("K32").Value = 3595.64
$ws.Range("M32").Value = -3308.64

$ws.Range("H45").Value = 3140.2666
$ws.Range("I45").Value = 2842.25
$ws.Range("K45").Value = 2842.25
$ws.Range("M45").Value = -2465.25

$ws.Range("H61").Value = 1057.2727
$ws.Range("I61").Value = 863
$ws.Range("K61").Value = 863
$ws.Range("M61").Value = -651

$ws.Range("H132").Value = 2130.85
$ws.Range("J132").Value = 4399.4614
$ws.Range("L132").Value = 13198.3842
$ws.Range("N132").Value = -18258.3842

$ws.Range("H136").Value = 1057.2727
$ws.Range("I136").Value = 863
$ws.Range("K136").Value = 2589
$ws.Range("M136").Value = -39

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7065.4165
$ws.Range("I20").Value = 6488.6
$ws.Range("K20").Value = 6488.6
$ws.Range("M20").Value = -6241.6

$ws.Range("H99").Value = 2336.6667
$ws.Range("I99").Value = 2005
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2005
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -507
$ws.Range("N99").Value = -5996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2004.4375
$ws.Range("I31").Value = 1619.9231
$ws.Range("K31").Value = 1619.9231
$ws.Range("M31").Value = -1324.9231

$ws.Range("H34").Value = 2004.4375
$ws.Range("I34").Value = 1619.9231
$ws.Range("K34").Value = 1619.9231
$ws.Range("M34").Value = -1417.9231

$ws.Range("H58").Value = 3006.75
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3006.75
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3006.75
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3412.75

$ws.Range("H105").Value = 2222
$ws.Range("I105").Value = 2222
$ws.Range("K105").Value = 2222
$ws.Range("M105").Value = -475

$ws.Range("H132").Value = 3229.8572
$ws.Range("I132").Value = 1903
$ws.Range("K132").Value = 5709
$ws.Range("M132").Value = -3179

$ws.Range("H134").Value = 3220.6667
$ws.Range("I134").Value = 3560.75
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 10682.25
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = -8147.25
$ws.Range("N134").Value = -6570

$ws.Range("H136").Value = 3006.75
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3006.75
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 9020.25
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -14120.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H113").Value = 1083.85
$ws.Range("J113").Value = 1037
$ws.Range("L113").Value = 3111
$ws.Range("N113").Value = -7451

$ws.Range("H121").Value = 526.8889
$ws.Range("I121").Value = 290.83334
$ws.Range("J121").Value = 999
$ws.Range("K121").Value = 872.5000200000001
$ws.Range("L121").Value = 2997
$ws.Range("M121").Value = 437.4999799999999
$ws.Range("N121").Value = -5617

$ws.Range("H122").Value = 101150
$ws.Range("J122").Value = 112288.445
$ws.Range("L122").Value = 1010596.005
$ws.Range("N122").Value = -1015496.005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14999
$ws.Range("I70").Value = 14999
$ws.Range("K70").Value = 14999
$ws.Range("M70").Value = -14729

$ws.Range("H73").Value = 14999
$ws.Range("I73").Value = 14999
$ws.Range("K73").Value = 14999
$ws.Range("M73").Value = -14063

$ws.Range("H97").Value = 1532.375
$ws.Range("I97").Value = 1692
$ws.Range("K97").Value = 1692
$ws.Range("M97").Value = -1196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4751.25
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 4005
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 4005
$ws.Range("M40").Value = -4864
$ws.Range("N40").Value = -4277

$ws.Range("H132").Value = 2534.5
$ws.Range("I132").Value = 1602.55
$ws.Range("K132").Value = 4807.65
$ws.Range("M132").Value = -2277.65

$ws.Range("H136").Value = 3812727.5
$ws.Range("I136").Value = 5336042.5
$ws.Range("J136").Value = 4440.8335
$ws.Range("K136").Value = 16008127.5
$ws.Range("L136").Value = 13322.5005
$ws.Range("M136").Value = -16005577.5
$ws.Range("N136").Value = -18422.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2216.2222
$ws.Range("I126").Value = 2076.3333
$ws.Range("J126").Value = 2496
$ws.Range("K126").Value = 6228.999899999999
$ws.Range("L126").Value = 7488
$ws.Range("M126").Value = -3758.999899999999
$ws.Range("N126").Value = -12428

$ws.Range("H132").Value = 1726.7097
$ws.Range("I132").Value = 1187.841
$ws.Range("J132").Value = 3043.9443
$ws.Range("K132").Value = 3563.523
$ws.Range("L132").Value = 9131.832900000001
$ws.Range("M132").Value = -1033.523
$ws.Range("N132").Value = -14191.8329

$ws.Range("H136").Value = 1466.8889
$ws.Range("I136").Value = 1119.5769
$ws.Range("J136").Value = 2369.9
$ws.Range("K136").Value = 3358.7307
$ws.Range("L136").Value = 7109.700000000001
$ws.Range("M136").Value = -808.7307000000001
$ws.Range("N136").Value = -12209.7
